# Fix ANV code and regenerate objects, fix mainrun automl code
#
# Adds a new "Feature group" style column (H/I/J) to the gcd_ml_features
# worksheet, tagging each feature row with its grp.* label. The
# "exclude" label that used to live in H10 is moved to J10 and replaced
# in H10 by the new "grp.anv." label.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New group labels for rows 2-5 (contact complementarity features)
$ws.Range("H2").Value = "grp.compl."
$ws.Range("H3").Value = "grp.compl."
$ws.Range("H4").Value = "grp.compl."
$ws.Range("H5").Value = "grp.compl."

# Kmer/triad related rows (7 and 8), entered before the GC.ij row (6)
# so the shared-string table picks up this exact ordering.
$ws.Range("H7").Value = "grp.kmer3."
$ws.Range("H8").Value = "grp.kmer1."

# Overall %GC row
$ws.Range("H6").Value = "grp.GC.ij"

# Per-region %GC columns: j first, then i
$ws.Range("I9").Value = "grp.GC.j_"
$ws.Range("H9").Value = "grp.GC.i_"

# ANV row (10): H10 used to hold "exclude" - move that to J10 and
# replace H10 with the new "grp.anv." label, carrying over the row's
# highlighted fill style.
$ws.Range("H10").Value = "grp.anv."
$ws.Range("J10").Value = "exclude"
$ws.Range("J10").Interior.Color = $ws.Range("H10").Interior.Color

# Restore the active selection to A9
$ws.Range("A9").Select()
